# Update NATMI LR-pairs sheet (Sema3c-Nrp2) with refreshed TPM figures.
# Sending-cluster "ECs" rows are gone; only FAPs/MuSCs remain as senders
# (6 data rows instead of 9), and every metric column reflects the new run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Sema3c"
$ws.Range("C2").Value = "Nrp2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 39.09670133333334
$ws.Range("H2").Value = 117.290104
$ws.Range("I2").Value = 0.9758026125363394
$ws.Range("J2").Value = 0.9758026125363395
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 18.43631966666667
$ws.Range("N2").Value = 55.308959
$ws.Range("O2").Value = 0.6034704469962782
$ws.Range("P2").Value = 0.6034704469962781
$ws.Range("Q2").Value = 720.7992836935263
$ws.Range("R2").Value = 6487.193553241737
$ws.Range("S2").Value = 0.5888680387674408
$ws.Range("T2").Value = 0.5888680387674408

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Sema3c"
$ws.Range("C3").Value = "Nrp2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 39.09670133333334
$ws.Range("H3").Value = 117.290104
$ws.Range("I3").Value = 0.9758026125363394
$ws.Range("J3").Value = 0.9758026125363395
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 5.347618333333333
$ws.Range("N3").Value = 16.042855
$ws.Range("O3").Value = 0.1750419652256785
$ws.Range("P3").Value = 0.1750419652256784
$ws.Range("Q3").Value = 209.0742368229911
$ws.Range("R3").Value = 1881.66813140692
$ws.Range("S3").Value = 0.1708064069707121
$ws.Range("T3").Value = 0.1708064069707121

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Sema3c"
$ws.Range("C4").Value = "Nrp2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 39.09670133333334
$ws.Range("H4").Value = 117.290104
$ws.Range("I4").Value = 0.9758026125363394
$ws.Range("J4").Value = 0.9758026125363395
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 6.766555
$ws.Range("N4").Value = 20.299665
$ws.Range("O4").Value = 0.2214875877780434
$ws.Range("P4").Value = 0.2214875877780434
$ws.Range("Q4").Value = 264.5499798905734
$ws.Range("R4").Value = 2380.94981901516
$ws.Range("S4").Value = 0.2161281667981865
$ws.Range("T4").Value = 0.2161281667981865

$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Sema3c"
$ws.Range("C5").Value = "Nrp2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.9694973333333333
$ws.Range("H5").Value = 2.908492
$ws.Range("I5").Value = 0.02419738746366056
$ws.Range("J5").Value = 0.02419738746366056
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 18.43631966666667
$ws.Range("N5").Value = 55.308959
$ws.Range("O5").Value = 0.6034704469962782
$ws.Range("P5").Value = 0.6034704469962781
$ws.Range("Q5").Value = 17.87396275331422
$ws.Range("R5").Value = 160.865664779828
$ws.Range("S5").Value = 0.01460240822883737
$ws.Range("T5").Value = 0.01460240822883737

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Sema3c"
$ws.Range("C6").Value = "Nrp2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.9694973333333333
$ws.Range("H6").Value = 2.908492
$ws.Range("I6").Value = 0.02419738746366056
$ws.Range("J6").Value = 0.02419738746366056
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 5.347618333333333
$ws.Range("N6").Value = 16.042855
$ws.Range("O6").Value = 0.1750419652256785
$ws.Range("P6").Value = 0.1750419652256784
$ws.Range("Q6").Value = 5.184501713851111
$ws.Range("R6").Value = 46.66051542466
$ws.Range("S6").Value = 0.004235558254966339
$ws.Range("T6").Value = 0.004235558254966338

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Sema3c"
$ws.Range("C7").Value = "Nrp2"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.9694973333333333
$ws.Range("H7").Value = 2.908492
$ws.Range("I7").Value = 0.02419738746366056
$ws.Range("J7").Value = 0.02419738746366056
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 6.766555
$ws.Range("N7").Value = 20.299665
$ws.Range("O7").Value = 0.2214875877780434
$ws.Range("P7").Value = 0.2214875877780434
$ws.Range("Q7").Value = 6.560157028353333
$ws.Range("R7").Value = 59.04141325518
$ws.Range("S7").Value = 0.005359420979856844
$ws.Range("T7").Value = 0.005359420979856843

# The old rows 8-10 (sending cluster "ECs") are removed entirely; the
# sheet now ends at row 7 (dimension A1:T7).
$ws.Rows("8:10").Delete()